# Agregué artículo sometido Andrés
# Insert a new supervision-table entry (Angela Rivero Valderrama & Sebastián
# Camilo Valenzuela) as the new second row of the "supervision" sheet,
# pushing all the existing entries down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record right below the header row.
$ws.Rows.Item(2).EntireRow.Insert()

$ws.Range("A2").Value = "BSc in Psychology"
$ws.Range("B2").Value = "2023 - 2024"
$ws.Range("C2").Value = "Angela Rivero Valderrama \& Sebastián Camilo Valenzuela"
$ws.Range("D2").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Range("E2").Value = "Research project: \textit{Preferencias por estímulos sexuales eróticos según género y la orientación sexual: un estudio con eye-tracking[Preferences for erotic sexual stimuli according to gender and sexual orientation: a study with eye-tracking]}"

# Match the formatting used by the rest of the table (wrapped, vertically
# centered text) and the row's natural wrap height.
$ws.Range("A2:E2").VerticalAlignment = -4108
$ws.Range("A2:E2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 43.2

$ws.Range("E2").Select()
$excel.ActiveWindow.ScrollColumn = 3
